$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data for rows 2-13 (columns A-T), reflecting the
# re-run of the NATMI pipeline with updated TPM values. Row layout:
# A=Sending cluster, B=Ligand symbol, C=Receptor symbol, D=Target cluster,
# E..T = numeric stats (expression / specificity columns).
$data = @(
    @("ECs", "Cntf", "Cntfr", "ECs", 2, 0.6666666666666666, 0.5906356666666667, 1.771907, 0.294823169192623, 0.294823169192623, 1, 0.3333333333333333, 0.03844333333333334, 0.11533, 0.006478890266597937, 0.006478890266597936, 0.02270600381222222, 0.20435403431, 0.001910126961249642, 0.001910126961249641),
    @("ECs", "Cntf", "Cntfr", "FAPs", 2, 0.6666666666666666, 0.5906356666666667, 1.771907, 0.294823169192623, 0.294823169192623, 3, 1, 5.666771333333333, 17.000314, 0.95502617622222, 0.9550261762222199, 3.346997264310889, 30.122975378798, 0.2815638439357473, 0.2815638439357473),
    @("ECs", "Cntf", "Cntfr", "MuSCs", 2, 0.6666666666666666, 0.5906356666666667, 1.771907, 0.294823169192623, 0.294823169192623, 3, 1, 0.2284146666666667, 0.685244, 0.03849493351118214, 0.03849493351118213, 0.1349098489231111, 1.214188640308, 0.01134919829562602, 0.01134919829562602),
    @("FAPs", "Cntf", "Cntfr", "ECs", 2, 0.6666666666666666, 0.5813396666666667, 1.744019, 0.2901829546991739, 0.2901829546991739, 1, 0.3333333333333333, 0.03844333333333334, 0.11533, 0.006478890266597937, 0.006478890266597936, 0.02234863458555556, 0.20113771127, 0.001880063520733108, 0.001880063520733108),
    @("FAPs", "Cntf", "Cntfr", "FAPs", 2, 0.6666666666666666, 0.5813396666666667, 1.744019, 0.2901829546991739, 0.2901829546991739, 3, 1, 5.666771333333333, 17.000314, 0.95502617622222, 0.9550261762222199, 3.294318957996222, 29.648870621966, 0.2771323176312177, 0.2771323176312177),
    @("FAPs", "Cntf", "Cntfr", "MuSCs", 2, 0.6666666666666666, 0.5813396666666667, 1.744019, 0.2901829546991739, 0.2901829546991739, 3, 1, 0.2284146666666667, 0.685244, 0.03849493351118214, 0.03849493351118213, 0.1327865061817778, 1.195078555636, 0.01117057354722308, 0.01117057354722308),
    @("MuSCs", "Cntf", "Cntfr", "ECs", 2, 0.6666666666666666, 0.817256, 2.451768, 0.4079435387325965, 0.4079435387325965, 1, 0.3333333333333333, 0.03844333333333334, 0.11533, 0.006478890266597937, 0.006478890266597936, 0.03141804482666667, 0.28276240344, 0.002643021422416138, 0.002643021422416137),
    @("MuSCs", "Cntf", "Cntfr", "FAPs", 2, 0.6666666666666666, 0.817256, 2.451768, 0.4079435387325965, 0.4079435387325965, 3, 1, 5.666771333333333, 17.000314, 0.95502617622222, 0.9550261762222199, 4.631202872794667, 41.680825855152, 0.3895967579103527, 0.3895967579103526),
    @("MuSCs", "Cntf", "Cntfr", "MuSCs", 2, 0.6666666666666666, 0.817256, 2.451768, 0.4079435387325965, 0.4079435387325965, 3, 1, 0.2284146666666667, 0.685244, 0.03849493351118214, 0.03849493351118213, 0.1866732568213333, 1.680059311392, 0.01570375939982766, 0.01570375939982765),
    @("Resolving-Mac", "Cntf", "Cntfr", "ECs", 2, 0.6666666666666666, 0.01412433333333333, 0.042373, 0.007050337375606629, 0.007050337375606628, 1, 0.3333333333333333, 0.03844333333333334, 0.11533, 0.006478890266597937, 0.006478890266597936, 0.0005429864544444445, 0.00488687809, [double]"4.567836219904942e-05", [double]"4.567836219904942e-05"),
    @("Resolving-Mac", "Cntf", "Cntfr", "FAPs", 2, 0.6666666666666666, 0.01412433333333333, 0.042373, 0.007050337375606629, 0.007050337375606628, 3, 1, 5.666771333333333, 17.000314, 0.95502617622222, 0.9550261762222199, 0.08003936723577779, 0.720354305122, 0.0067332567449022, 0.006733256744902198),
    @("Resolving-Mac", "Cntf", "Cntfr", "MuSCs", 2, 0.6666666666666666, 0.01412433333333333, 0.042373, 0.007050337375606629, 0.007050337375606628, 3, 1, 0.2284146666666667, 0.685244, 0.03849493351118214, 0.03849493351118213, 0.003226204890222222, 0.029035844012, 0.0002714022685053796, 0.0002714022685053795)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $startRow + $i
    for ($c = 1; $c -le $row.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
}
